$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Temperature-Rearing_score" column (P).
# This shifts the old column P (and all its data) one column to the right, into Q.
$ws.Columns("P:P").Insert()

# Give the newly inserted column P its header.
$ws.Range("P1").Value = "Temperature-AdultHolding_score"

# Update the new "Temperature-AdultHolding_score" (P) and
# "Temperature-Rearing_score" (Q, shifted from old P) columns with the
# corrected values for each reach row.
$ws.Range("P2").Value = "Inf"
$ws.Range("Q2").Value = 5

$ws.Range("P3").Value = 5
$ws.Range("Q3").Value = "Inf"

$ws.Range("P4").Value = "Inf"
$ws.Range("Q4").Value = 5

$ws.Range("P5").Value = 5
$ws.Range("Q5").Value = "Inf"
